# Implemented further mutations of articles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New content: add the new mutation rows/cells first, in the order the
# --- new shared strings need to appear (59..64) ---
$ws.Range("B14").Value = "voor onderdeel A"
$ws.Range("D13").Value = "Artikel 44"
$ws.Range("D14").Value = "Artikel 8.14a"
$ws.Range("D15").Value = "artikel 8.9"
$ws.Range("D16").Value = "artikel 9a"
$ws.Range("A18").ClearContents()
$ws.Range("C20").Value = "s"

# --- Highlight fills (theme colors) on the "Article (roman)" and
# --- "Article (numeric)" columns to mark resolved / in-progress mutations ---
# Green = xlThemeColorAccent6 (theme index 9 in the XML clrScheme)
$greenCells = @("A5","D5","A6","A7","A8","D8","A9","D9","D10","D11","D13","D14","D15","D16")
foreach ($addr in $greenCells) {
    $ws.Range($addr).Interior.ThemeColor = 10
}

# Orange = xlThemeColorAccent2 (theme index 5 in the XML clrScheme)
$orangeCells = @("D6","D7","D12")
foreach ($addr in $orangeCells) {
    $ws.Range($addr).Interior.ThemeColor = 6
}

# --- Selection moved to C17 ---
$ws.Range("C17").Select() | Out-Null

Write-Host "Applied tweede kamer mutation updates"
